$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 7).Value = 3.2
$ws.Cells.Item(3, 9).Value = 2.4
$ws.Cells.Item(3, 10).Value = 1.1
$ws.Cells.Item(3, 11).Value = 7
$ws.Cells.Item(3, 20).Value = 7.5
$ws.Cells.Item(3, 22).Value = 12

# Row 7
$ws.Cells.Item(7, 7).Value = 1.47
$ws.Cells.Item(7, 8).Value = 3.95
$ws.Cells.Item(7, 9).Value = 6.2
$ws.Cells.Item(7, 18).Value = 2.12
$ws.Cells.Item(7, 19).Value = 1.57
$ws.Cells.Item(7, 20).Value = 5.6
$ws.Cells.Item(7, 21).Value = 6.1
$ws.Cells.Item(7, 22).Value = 8.5
$ws.Cells.Item(7, 23).Value = 9.5
$ws.Cells.Item(7, 27).Value = 7.9
$ws.Cells.Item(7, 28).Value = 23
$ws.Cells.Item(7, 31).Value = 13.5
$ws.Cells.Item(7, 32).Value = 37
$ws.Cells.Item(7, 33).Value = 21
$ws.Cells.Item(7, 34).Value = 150
$ws.Cells.Item(7, 35).Value = 80
$ws.Cells.Item(7, 36).Value = 90

# Row 8
$ws.Cells.Item(8, 8).Value = 3.65
$ws.Cells.Item(8, 11).Value = 6.6
$ws.Cells.Item(8, 12).Value = 1.36
$ws.Cells.Item(8, 13).Value = 2.92
$ws.Cells.Item(8, 14).Value = 2.05
$ws.Cells.Item(8, 15).Value = 1.7
$ws.Cells.Item(8, 16).Value = 1.42
$ws.Cells.Item(8, 17).Value = 2.65
$ws.Cells.Item(8, 18).Value = 2.12
$ws.Cells.Item(8, 19).Value = 1.65
$ws.Cells.Item(8, 20).Value = 14.5
$ws.Cells.Item(8, 21).Value = 40
$ws.Cells.Item(8, 22).Value = 21
$ws.Cells.Item(8, 23).Value = 150
$ws.Cells.Item(8, 24).Value = 80
$ws.Cells.Item(8, 25).Value = 80
$ws.Cells.Item(8, 26).Value = 6.6
$ws.Cells.Item(8, 28).Value = 20
$ws.Cells.Item(8, 32).Value = 6.2
$ws.Cells.Item(8, 33).Value = 8.25
$ws.Cells.Item(8, 34).Value = 10.25

# Row 10
$ws.Cells.Item(10, 10).Value = 1.04
$ws.Cells.Item(10, 12).Value = 1.25
$ws.Cells.Item(10, 18).Value = 1.75

# Row 13
$ws.Cells.Item(13, 10).Value = 1.01
$ws.Cells.Item(13, 11).Value = 23
$ws.Cells.Item(13, 12).Value = 1.11
$ws.Cells.Item(13, 13).Value = 6.5
$ws.Cells.Item(13, 14).Value = 1.4
$ws.Cells.Item(13, 15).Value = 2.88

# Row 16
$ws.Cells.Item(16, 7).Value = 2.82
$ws.Cells.Item(16, 8).Value = 2.87
$ws.Cells.Item(16, 10).Value = 1.11
$ws.Cells.Item(16, 11).Value = 5.8
$ws.Cells.Item(16, 13).Value = 2.5
$ws.Cells.Item(16, 14).Value = 2.42
$ws.Cells.Item(16, 15).Value = 1.52
$ws.Cells.Item(16, 16).Value = 1.53
$ws.Cells.Item(16, 17).Value = 2.4
$ws.Cells.Item(16, 18).Value = 1.98
$ws.Cells.Item(16, 19).Value = 1.75
$ws.Cells.Item(16, 20).Value = 7.1
$ws.Cells.Item(16, 22).Value = 11.25
$ws.Cells.Item(16, 23).Value = 37
$ws.Cells.Item(16, 24).Value = 30
$ws.Cells.Item(16, 26).Value = 5.8
$ws.Cells.Item(16, 27).Value = 5.9
$ws.Cells.Item(16, 30).Value = 1250
$ws.Cells.Item(16, 34).Value = 35

# Row 20
$ws.Cells.Item(20, 10).Value = 1.01
$ws.Cells.Item(20, 11).Value = 34
$ws.Cells.Item(20, 12).Value = 1.03
$ws.Cells.Item(20, 13).Value = 17
$ws.Cells.Item(20, 14).Value = 1.13
$ws.Cells.Item(20, 15).Value = 6
$ws.Cells.Item(20, 18).Value = 1.83
$ws.Cells.Item(20, 19).Value = 1.83
$ws.Cells.Item(20, 21).Value = 12
$ws.Cells.Item(20, 23).Value = 9
$ws.Cells.Item(20, 27).Value = 29
$ws.Cells.Item(20, 30).Value = 201
$ws.Cells.Item(20, 36).Value = 81

# Row 22
$ws.Cells.Item(22, 8).Value = 8.5
$ws.Cells.Item(22, 9).Value = 26
$ws.Cells.Item(22, 13).Value = 7.1
$ws.Cells.Item(22, 14).Value = 1.23
$ws.Cells.Item(22, 15).Value = 3.75
$ws.Cells.Item(22, 17).Value = 4.7
$ws.Cells.Item(22, 18).Value = 2.1
$ws.Cells.Item(22, 19).Value = 1.65
$ws.Cells.Item(22, 20).Value = 13.5
$ws.Cells.Item(22, 21).Value = 8
$ws.Cells.Item(22, 22).Value = 13
$ws.Cells.Item(22, 23).Value = 7
$ws.Cells.Item(22, 24).Value = 11
$ws.Cells.Item(22, 25).Value = 32
$ws.Cells.Item(22, 27).Value = 23
$ws.Cells.Item(22, 28).Value = 40
$ws.Cells.Item(22, 29).Value = 150
$ws.Cells.Item(22, 30).Value = 800
$ws.Cells.Item(22, 32).Value = 450
$ws.Cells.Item(22, 35).Value = 500
$ws.Cells.Item(22, 36).Value = 200

# Row 23
$ws.Cells.Item(23, 7).Value = 3.65
$ws.Cells.Item(23, 10).Value = 1.04
$ws.Cells.Item(23, 11).Value = 8.75
$ws.Cells.Item(23, 12).Value = 1.19
$ws.Cells.Item(23, 13).Value = 4.2
$ws.Cells.Item(23, 14).Value = 1.57
$ws.Cells.Item(23, 15).Value = 2.25
$ws.Cells.Item(23, 18).Value = 1.52
$ws.Cells.Item(23, 19).Value = 2.37
$ws.Cells.Item(23, 21).Value = 24
$ws.Cells.Item(23, 26).Value = 8.75
$ws.Cells.Item(23, 27).Value = 7.4
$ws.Cells.Item(23, 30).Value = 250
$ws.Cells.Item(23, 31).Value = 9.75
$ws.Cells.Item(23, 32).Value = 11
$ws.Cells.Item(23, 35).Value = 13.5
$ws.Cells.Item(23, 36).Value = 19

# Row 24
$ws.Cells.Item(24, 7).Value = 1.57
$ws.Cells.Item(24, 8).Value = 4.05
$ws.Cells.Item(24, 9).Value = 5
$ws.Cells.Item(24, 11).Value = 9.5
$ws.Cells.Item(24, 12).Value = 1.14
$ws.Cells.Item(24, 13).Value = 5
$ws.Cells.Item(24, 14).Value = 1.44
$ws.Cells.Item(24, 15).Value = 2.62
$ws.Cells.Item(24, 16).Value = 1.25
$ws.Cells.Item(24, 17).Value = 3.55
$ws.Cells.Item(24, 18).Value = 1.47
$ws.Cells.Item(24, 19).Value = 2.5
$ws.Cells.Item(24, 20).Value = 11.25
$ws.Cells.Item(24, 21).Value = 10.5
$ws.Cells.Item(24, 22).Value = 8.25
$ws.Cells.Item(24, 23).Value = 14
$ws.Cells.Item(24, 24).Value = 11
$ws.Cells.Item(24, 25).Value = 16.5
$ws.Cells.Item(24, 26).Value = 9.5
$ws.Cells.Item(24, 27).Value = 8.75
$ws.Cells.Item(24, 28).Value = 12.5
$ws.Cells.Item(24, 29).Value = 37
$ws.Cells.Item(24, 30).Value = 200
$ws.Cells.Item(24, 31).Value = 20
$ws.Cells.Item(24, 32).Value = 37
$ws.Cells.Item(24, 33).Value = 16
$ws.Cells.Item(24, 34).Value = 100

# Row 25
$ws.Cells.Item(25, 7).Value = 4.35
$ws.Cells.Item(25, 8).Value = 3.6
$ws.Cells.Item(25, 9).Value = 1.75
$ws.Cells.Item(25, 10).Value = 1.04
$ws.Cells.Item(25, 11).Value = 8.75
$ws.Cells.Item(25, 12).Value = 1.19
$ws.Cells.Item(25, 13).Value = 4.25
$ws.Cells.Item(25, 14).Value = 1.57
$ws.Cells.Item(25, 15).Value = 2.25
$ws.Cells.Item(25, 16).Value = 1.32
$ws.Cells.Item(25, 17).Value = 3.1
$ws.Cells.Item(25, 18).Value = 1.53
$ws.Cells.Item(25, 19).Value = 2.32
$ws.Cells.Item(25, 20).Value = 17
$ws.Cells.Item(25, 22).Value = 14
$ws.Cells.Item(25, 26).Value = 8.75
$ws.Cells.Item(25, 27).Value = 7.4
$ws.Cells.Item(25, 28).Value = 12
$ws.Cells.Item(25, 29).Value = 40
$ws.Cells.Item(25, 30).Value = 250
$ws.Cells.Item(25, 31).Value = 9.25
$ws.Cells.Item(25, 32).Value = 10
